# Generate Report for Handback
#
# Marks a.md / b.md as handed back (status text flips from "Ready for
# handoff" to "Handed back: in sync with en-US" in the Overview sheet and
# in both language sheets), and fills in the "Latest Target File",
# "Latest Handback File" and "Latest Handback DateTime" columns on the
# zh-cn / de-de sheets for those two rows.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---- Overview sheet: update Status column for a.md / b.md rows ----
$ovw = $wb.Worksheets.Item("Overview")
$ovw.Range("B2").Value = $newStatus
$ovw.Range("C2").Value = $newStatus
$ovw.Range("B3").Value = $newStatus
$ovw.Range("C3").Value = $newStatus

# ---- zh-cn sheet ----
$zh = $wb.Worksheets.Item("zh-cn")
$zhXlf = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/865a08e6aa25208950fcfe2e04d87511d0e484b2/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$mdTarget = "https://github.com/OpenLocalizationTest/oltest/blob/65f7ce8b99da9587251708640009eef78f7bf43b/e2e/a.md"

foreach ($row in 2, 3) {
    $zh.Cells.Item($row, 2).Value = $newStatus

    # Latest Target File (E)
    $zh.Hyperlinks.Add($zh.Cells.Item($row, 5), $mdTarget, $null, $null, "a.md")
    # Latest Handback File (F)
    $zh.Hyperlinks.Add($zh.Cells.Item($row, 6), $zhXlf, $null, $null, "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf")
    # Latest Handback DateTime (G)
    $zh.Cells.Item($row, 7).Value = "2016-02-18 09:45:38"
}

# ---- de-de sheet ----
$de = $wb.Worksheets.Item("de-de")
$deXlf = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/56fc8672983f6421f83125d23d7fd254c300d320/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

foreach ($row in 2, 3) {
    $de.Cells.Item($row, 2).Value = $newStatus

    # Latest Target File (E)
    $de.Hyperlinks.Add($de.Cells.Item($row, 5), $mdTarget, $null, $null, "a.md")
    # Latest Handback File (F)
    $de.Hyperlinks.Add($de.Cells.Item($row, 6), $deXlf, $null, $null, "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf")
    # Latest Handback DateTime (G)
    $de.Cells.Item($row, 7).Value = "2016-02-18 09:46:02"
}

$wb.Save()
